$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 99.15000000000001
$ws.Range("I5").Value = 98.375
$ws.Range("K5").Value = 98.375
$ws.Range("M5").Value = 16.625
# Row 116
$ws.Range("H116").Value = 3353.318
$ws.Range("I116").Value = 3239.5881
$ws.Range("K116").Value = 3239.5881
$ws.Range("M116").Value = 202.4119000000001
# Row 134
$ws.Range("H134").Value = 82210.91
$ws.Range("J134").Value = 82210.91
$ws.Range("L134").Value = 82210.91
$ws.Range("N134").Value = -92350.91
# Row 138
$ws.Range("H138").Value = 2475.59
$ws.Range("I138").Value = 904.8
$ws.Range("J138").Value = 2999.1868
$ws.Range("K138").Value = 2714.4
$ws.Range("L138").Value = 8997.5604
$ws.Range("M138").Value = 2425.6
$ws.Range("N138").Value = -19277.5604

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 33653.03
$ws.Range("I2").Value = 46814
$ws.Range("J2").Value = 1481.7778
$ws.Range("K2").Value = 46814
$ws.Range("L2").Value = 1481.7778
$ws.Range("M2").Value = -46701
$ws.Range("N2").Value = -1707.7778
# Row 32
$ws.Range("H32").Value = 2130.05
$ws.Range("I32").Value = 1994.2235
$ws.Range("J32").Value = 2899.7334
$ws.Range("K32").Value = 1994.2235
$ws.Range("L32").Value = 2899.7334
$ws.Range("M32").Value = -1707.2235
$ws.Range("N32").Value = -3473.7334
# Row 110
$ws.Range("H110").Value = 667.5833
$ws.Range("I110").Value = 656.7778
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 656.7778
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1388.2222
$ws.Range("N110").Value = -4790
# Row 116
$ws.Range("H116").Value = 33653.03
$ws.Range("I116").Value = 46814
$ws.Range("J116").Value = 1481.7778
$ws.Range("K116").Value = 46814
$ws.Range("L116").Value = 1481.7778
$ws.Range("M116").Value = -44520
$ws.Range("N116").Value = -6069.7778
# Row 122
$ws.Range("H122").Value = 690.4
$ws.Range("I122").Value = 690.4
$ws.Range("K122").Value = 2071.2
$ws.Range("M122").Value = 378.8000000000002
# Row 132
$ws.Range("H132").Value = 1255.9318
$ws.Range("I132").Value = 755.91174
$ws.Range("J132").Value = 2956
$ws.Range("K132").Value = 2267.73522
$ws.Range("L132").Value = 8868
$ws.Range("M132").Value = 262.26478
$ws.Range("N132").Value = -13928

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 33653.03
$ws.Range("I3").Value = 46814
$ws.Range("J3").Value = 1481.7778
$ws.Range("K3").Value = 46814
$ws.Range("L3").Value = 1481.7778
$ws.Range("M3").Value = -46700
$ws.Range("N3").Value = -1709.7778
# Row 134
$ws.Range("H134").Value = 20045.629
$ws.Range("I134").Value = 1537.1621
$ws.Range("K134").Value = 4611.4863
$ws.Range("M134").Value = -2076.4863

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 2781.1177
$ws.Range("I86").Value = 2177.9
$ws.Range("J86").Value = 3642.8572
$ws.Range("K86").Value = 2177.9
$ws.Range("L86").Value = 3642.8572
$ws.Range("M86").Value = -1054.9
$ws.Range("N86").Value = -5888.8572
# Row 89
$ws.Range("H89").Value = 2781.1177
$ws.Range("I89").Value = 2177.9
$ws.Range("J89").Value = 3642.8572
$ws.Range("K89").Value = 10889.5
$ws.Range("L89").Value = 18214.286
$ws.Range("M89").Value = -5273.5
$ws.Range("N89").Value = -29446.286
# Row 122
$ws.Range("H122").Value = 1365.7273
$ws.Range("I122").Value = 1210.04
$ws.Range("J122").Value = 1852.25
$ws.Range("K122").Value = 3630.12
$ws.Range("L122").Value = 5556.75
$ws.Range("M122").Value = -1180.12
$ws.Range("N122").Value = -10456.75
# Row 132
$ws.Range("H132").Value = 1543.9807
$ws.Range("I132").Value = 799.0714
$ws.Range("J132").Value = 4672.6
$ws.Range("K132").Value = 2397.2142
$ws.Range("L132").Value = 14017.8
$ws.Range("M132").Value = 132.7857999999997
$ws.Range("N132").Value = -19077.8

$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 72230.87
$ws.Range("I137").Value = 2443
$ws.Range("J137").Value = 211806.6
$ws.Range("K137").Value = 7329
$ws.Range("L137").Value = 635419.8
$ws.Range("M137").Value = -2229
$ws.Range("N137").Value = -645619.8

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3995.6667
$ws.Range("I80").Value = 3776.9565
$ws.Range("J80").Value = 4714.2856
$ws.Range("K80").Value = 3776.9565
$ws.Range("L80").Value = 4714.2856
$ws.Range("M80").Value = -2778.9565
$ws.Range("N80").Value = -6710.2856
# Row 83
$ws.Range("H83").Value = 3995.6667
$ws.Range("I83").Value = 3776.9565
$ws.Range("J83").Value = 4714.2856
$ws.Range("K83").Value = 18884.7825
$ws.Range("L83").Value = 23571.428
$ws.Range("M83").Value = -13892.7825
$ws.Range("N83").Value = -33555.428
# Row 122
$ws.Range("H122").Value = 5901978
$ws.Range("I122").Value = 5702519
$ws.Range("J122").Value = 6251031
$ws.Range("K122").Value = 17107557
$ws.Range("L122").Value = 18753093
$ws.Range("M122").Value = -17105107
$ws.Range("N122").Value = -18757993
# Row 123
$ws.Range("H123").Value = 11162.75
$ws.Range("J123").Value = 11162.75
$ws.Range("L123").Value = 11162.75
$ws.Range("N123").Value = -16062.75
# Row 132
$ws.Range("H132").Value = 1929.1892
$ws.Range("I132").Value = 1715.0625
$ws.Range("J132").Value = 3299.6
$ws.Range("K132").Value = 5145.1875
$ws.Range("L132").Value = 9898.799999999999
$ws.Range("M132").Value = -2615.1875
$ws.Range("N132").Value = -14958.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5293005
$ws.Range("I7").Value = 2084.9285
$ws.Range("J7").Value = 15874845
$ws.Range("K7").Value = 2084.9285
$ws.Range("L7").Value = 15874845
$ws.Range("M7").Value = -1972.9285
$ws.Range("N7").Value = -15875069
# Row 61
$ws.Range("H61").Value = 5910.273
$ws.Range("I61").Value = 7757.875
$ws.Range("K61").Value = 7757.875
$ws.Range("M61").Value = -7555.875
# Row 113
$ws.Range("H113").Value = 5910.273
$ws.Range("I113").Value = 7757.875
$ws.Range("K113").Value = 7757.875
$ws.Range("M113").Value = -5587.875
# Row 122
$ws.Range("H122").Value = 11699
$ws.Range("I122").Value = 15828.571
$ws.Range("J122").Value = 2063.3333
$ws.Range("K122").Value = 47485.713
$ws.Range("L122").Value = 6189.999899999999
$ws.Range("M122").Value = -45035.713
$ws.Range("N122").Value = -11089.9999
# Row 126
$ws.Range("H126").Value = 5293005
$ws.Range("I126").Value = 2084.9285
$ws.Range("J126").Value = 15874845
$ws.Range("K126").Value = 6254.7855
$ws.Range("L126").Value = 47624535
$ws.Range("M126").Value = -3784.7855
$ws.Range("N126").Value = -47629475
# Row 132
$ws.Range("H132").Value = 2553.913
$ws.Range("I132").Value = 2270.6316
$ws.Range("J132").Value = 3899.5
$ws.Range("K132").Value = 6811.8948
$ws.Range("L132").Value = 11698.5
$ws.Range("M132").Value = -4281.8948
$ws.Range("N132").Value = -16758.5
# Row 136
$ws.Range("H136").Value = 3209.524
$ws.Range("I136").Value = 1482.6471
$ws.Range("J136").Value = 10548.75
$ws.Range("K136").Value = 4447.9413
$ws.Range("L136").Value = 31646.25
$ws.Range("M136").Value = -1897.9413
$ws.Range("N136").Value = -36746.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1355.6923
$ws.Range("I122").Value = 1352.1818
$ws.Range("J122").Value = 1375
$ws.Range("K122").Value = 4056.5454
$ws.Range("L122").Value = 4125
$ws.Range("M122").Value = -1606.5454
$ws.Range("N122").Value = -9025
# Row 132
$ws.Range("H132").Value = 1023.21124
$ws.Range("I132").Value = 737.3333
$ws.Range("J132").Value = 1437.2413
$ws.Range("K132").Value = 2211.9999
$ws.Range("L132").Value = 4311.7239
$ws.Range("M132").Value = 318.0001000000002
$ws.Range("N132").Value = -9371.723900000001
# Row 136
$ws.Range("H136").Value = 1326.2572
$ws.Range("I136").Value = 1263.1017
$ws.Range("J136").Value = 1665
$ws.Range("K136").Value = 3789.3051
$ws.Range("L136").Value = 4995
$ws.Range("M136").Value = -1239.3051
$ws.Range("N136").Value = -10095
